$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Insert 3 new blank rows above the old row 120 block -----------------
# (old rows 120/121/122 shift down to 123/124/125; this also pushes the
#  per-column formatting of row 119 into the new blank rows, matching
#  Excel's "insert copies format from the row above" behaviour)
$ws.Rows.Item(120).Resize(3).Insert(-4121)

# --- Fill the newly-typed data row (row 119) ------------------------------
$ws.Cells.Item(119,1).Value2 = "CW3M C751"
$ws.Cells.Item(119,2).Value2 = "Demo_Baseline 2010-18"
$ws.Cells.Item(119,2).WrapText = $true
$ws.Cells.Item(119,3).Value2 = "2010-18"

$ws.Cells.Item(119,4).NumberFormat  = "0.00"
$ws.Cells.Item(119,4).Value2  = 936.57170944444442
$ws.Cells.Item(119,5).NumberFormat  = "0.00"
$ws.Cells.Item(119,5).Value2  = 1890.2624918888889
$ws.Cells.Item(119,6).NumberFormat  = "0.00"
$ws.Cells.Item(119,6).Value2  = 0.58615577777777772
$ws.Cells.Item(119,7).NumberFormat  = "0.00"
$ws.Cells.Item(119,7).Value2  = 270.41205844444437
$ws.Cells.Item(119,8).NumberFormat  = "0.00"
$ws.Cells.Item(119,8).Value2  = 0
$ws.Cells.Item(119,9).NumberFormat  = "0.00"
$ws.Cells.Item(119,9).Value2  = 0.20794977777777776
$ws.Cells.Item(119,10).NumberFormat = "0.00"
$ws.Cells.Item(119,10).Value2 = 8.1971030000000003
$ws.Cells.Item(119,11).NumberFormat = "0.00"
$ws.Cells.Item(119,11).Value2 = 662.99974244444445
$ws.Cells.Item(119,12).NumberFormat = "0.00"
$ws.Cells.Item(119,12).Value2 = 80.184555777777788
$ws.Cells.Item(119,13).NumberFormat = "0.00"
$ws.Cells.Item(119,13).Value2 = 1407.6443413333334
$ws.Cells.Item(119,14).NumberFormat = "0.00"
$ws.Cells.Item(119,14).Value2 = 938.96002866666663
$ws.Cells.Item(119,15).NumberFormat = "0"
$ws.Cells.Item(119,15).Value2 = 3872.6727430000005
$ws.Cells.Item(119,16).NumberFormat = "0"
$ws.Cells.Item(119,16).Value2 = 27393.764540111111
$ws.Cells.Item(119,17).NumberFormat = "0.00"
$ws.Cells.Item(119,17).Value2 = -0.054594222222222016
$ws.Cells.Item(119,18).NumberFormat = "0.000000"
$ws.Cells.Item(119,18).Value2 = -0.000023555555555555631

# --- Apply the same per-column number formats to the 3 blank rows below --
$ws.Range("D120:N122").NumberFormat = "0.00"
$ws.Range("O120:P122").NumberFormat = "0"
$ws.Range("Q120:Q122").NumberFormat = "0.00"
$ws.Range("R120:R122").NumberFormat = "0.000000"

# --- Match the author's final selection -----------------------------------
$ws.Range("B119").Select()
